$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values (e.g. "1.00", "0.0527")
# keep their exact original text representation instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.578.00"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "2.381.87"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "552.75"
$ws.Range("E5").Value = "  +2.12%  "

$ws.Range("D6").Value = "140.98"
$ws.Range("E6").Value = "  +3.27%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("D9").Value = "2.385.01"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +4.33%  "

$ws.Range("E11").Value = "  +2.09%  "

$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  +2.49%  "

$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +3.73%  "

$ws.Range("D14").Value = "25.72"
$ws.Range("E14").Value = "  +5.48%  "

$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +7.92%  "

$ws.Range("D16").Value = "2.813.70"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").Value = "61.534.71"
$ws.Range("E17").Value = "  +2.22%  "

$ws.Range("D18").Value = "2.383.34"
$ws.Range("E18").Value = "  +1.56%  "

$ws.Range("D19").Value = "10.99"
$ws.Range("E19").Value = "  +4.28%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "322.87"
$ws.Range("E20").Value = "  +3.46%  "

$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "4.17"
$ws.Range("E21").Value = "  +2.74%  "

$ws.Range("D22").Value = "6.69"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -5.27%  "

$ws.Range("D25").Value = "64.29"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").Value = "8.85"
$ws.Range("E26").Value = "  +3.20%  "

$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "2.500.20"
$ws.Range("E28").Value = "  +1.49%  "

$ws.Range("D29").Value = "533.18"
$ws.Range("E29").Value = "  +7.34%  "

$ws.Range("D30").Value = "8.26"
$ws.Range("E30").Value = "  +4.55%  "

$ws.Range("D31").Value = "0.0₃0914"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("E34").Value = "  +4.05%  "

$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "5.72"
$ws.Range("E36").Value = "  +9.81%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "4.75"
$ws.Range("E38").Value = "  +4.20%  "

$ws.Range("E39").Value = "  +8.18%  "

$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("D41").Value = "18.57"
$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("D42").Value = "146.39"
$ws.Range("E42").Value = "  +6.77%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "41.47"
$ws.Range("E44").Value = "  +4.08%  "

$ws.Range("D45").Value = "149.34"
$ws.Range("E45").Value = "  +5.70%  "

$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("D47").Value = "3.60"
$ws.Range("E47").Value = "  +2.45%  "

$ws.Range("D48").Value = "0.0527"
$ws.Range("E48").Value = "  +4.01%  "

$ws.Range("D49").Value = "20.02"
$ws.Range("E49").Value = "  +3.61%  "

$ws.Range("D50").Value = "0.585"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("D51").Value = "0.0906"
$ws.Range("E51").Value = "  +1.66%  "

